$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I3").Value = -0.195105760602844
$ws.Range("J3").Value = 0.6081287845824713
$ws.Range("K3").Value = 0.451901484868166
$ws.Range("L3").Value = 2.517552118679247

$ws.Range("I20").Value = 0.035706369411352
$ws.Range("J20").Value = 0.5931560001382098
$ws.Range("K20").Value = 0.1487503366435251
$ws.Range("L20").Value = 2.435309308996479
